# Apply the recorded changes to the "Artfynd" worksheet.
# The edit re-shuffles several existing observation rows (their species /
# coordinate / date data moves between rows while the location columns
# P,S,T,U,V,W and the reporter columns AW,AX stay put), bumps a single
# "Taxonsorteringsordning" value in row 12, and appends a brand new
# observation as row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng) {
    # Force text (not number/date) interpretation for the given single cell.
    $rng.NumberFormat = "@"
}

# ---------------------------------------------------------------------
# Row 2  (was row 3's data)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 111780621
$ws.Range("B2").Value = 56543
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 103021
$ws.Range("F2").Value = "Talltita"
$ws.Range("G2").Value = "Poecile montanus"
$ws.Range("H2").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q2").Value = 707631
$ws.Range("R2").Value = 7397278
Set-TextCell $ws.Range("Y2")
$ws.Range("Y2").Value = "2023-08-29"
Set-TextCell $ws.Range("AA2")
$ws.Range("AA2").Value = "2023-08-29"

# ---------------------------------------------------------------------
# Row 3  (was row 2's data)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 111780628
$ws.Range("B3").Value = 78604
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 6461
$ws.Range("F3").Value = "Norrlandslav"
$ws.Range("G3").Value = "Nephroma arcticum"
$ws.Range("H3").Value = "(L.) Torss."
$ws.Range("Q3").Value = 707614
$ws.Range("R3").Value = 7397255
Set-TextCell $ws.Range("Y3")
$ws.Range("Y3").Value = "2023-08-29"
Set-TextCell $ws.Range("AA3")
$ws.Range("AA3").Value = "2023-08-29"

# ---------------------------------------------------------------------
# Row 6  (was row 8's data)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 111816145
$ws.Range("B6").Value = 77597
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 864
$ws.Range("F6").Value = "Knottrig blåslav"
$ws.Range("G6").Value = "Hypogymnia bitteri"
$ws.Range("H6").Value = "(Lynge) Ahti"
$ws.Range("Q6").Value = 707627
$ws.Range("R6").Value = 7397312
Set-TextCell $ws.Range("Y6")
$ws.Range("Y6").Value = "2023-08-22"
Set-TextCell $ws.Range("AA6")
$ws.Range("AA6").Value = "2023-08-22"

# ---------------------------------------------------------------------
# Row 8  (was row 10's data)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 111816137
$ws.Range("B8").Value = 90658
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 4361
$ws.Range("F8").Value = "Orange taggsvamp"
$ws.Range("G8").Value = "Hydnellum aurantiacum"
$ws.Range("H8").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q8").Value = 707609
$ws.Range("R8").Value = 7397264
Set-TextCell $ws.Range("Y8")
$ws.Range("Y8").Value = "2023-08-22"
Set-TextCell $ws.Range("AA8")
$ws.Range("AA8").Value = "2023-08-22"

# ---------------------------------------------------------------------
# Row 9  (was row 11's data)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 111816119
$ws.Range("B9").Value = 56543
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 103021
$ws.Range("F9").Value = "Talltita"
$ws.Range("G9").Value = "Poecile montanus"
$ws.Range("H9").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q9").Value = 707596
$ws.Range("R9").Value = 7397263
Set-TextCell $ws.Range("Y9")
$ws.Range("Y9").Value = "2023-08-22"
Set-TextCell $ws.Range("AA9")
$ws.Range("AA9").Value = "2023-08-22"

# ---------------------------------------------------------------------
# Row 10 (was row 9's data)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 111816142
$ws.Range("B10").Value = 78604
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 6461
$ws.Range("F10").Value = "Norrlandslav"
$ws.Range("G10").Value = "Nephroma arcticum"
$ws.Range("H10").Value = "(L.) Torss."
$ws.Range("Q10").Value = 707613
$ws.Range("R10").Value = 7397270
Set-TextCell $ws.Range("Y10")
$ws.Range("Y10").Value = "2023-08-29"
Set-TextCell $ws.Range("AA10")
$ws.Range("AA10").Value = "2023-08-29"

# ---------------------------------------------------------------------
# Row 11 (was row 6's data)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 111816118
$ws.Range("B11").Value = 78107
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6453
$ws.Range("F11").Value = "Vedskivlav"
$ws.Range("G11").Value = "Hertelidea botryosa"
$ws.Range("H11").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q11").Value = 707670
$ws.Range("R11").Value = 7397328
Set-TextCell $ws.Range("Y11")
$ws.Range("Y11").Value = "2023-08-22"
Set-TextCell $ws.Range("AA11")
$ws.Range("AA11").Value = "2023-08-22"

# The (blank) "Bestämningsmetod" marker cell rides along with the moved
# species data for rows 8/9 -> it disappears from its old spot once the
# row no longer holds that particular record.
$ws.Range("AF8").ClearContents()
$ws.Range("AF9").ClearContents()

# ---------------------------------------------------------------------
# Row 12: only the "Taxonsorteringsordning" (column B) value changes.
# ---------------------------------------------------------------------
$ws.Range("B12").Value = 55643

# ---------------------------------------------------------------------
# Row 13: brand new observation appended at the bottom of the table.
# ---------------------------------------------------------------------
$ws.Range("A13").Value = 112432492
$ws.Range("B13").Value = 90823
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5966
$ws.Range("F13").Value = "Motaggsvamp"
$ws.Range("G13").Value = "Sarcodon squamosus"
$ws.Range("H13").Value = "(Schaeff.) Quél."
$ws.Range("P13").Value = "Saskam, Lu lm"
$ws.Range("Q13").Value = 707528
$ws.Range("R13").Value = 7397284
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Norrbotten"
$ws.Range("U13").Value = "Jokkmokk"
$ws.Range("V13").Value = "Lule lappmark"
$ws.Range("W13").Value = "Jokkmokk"
Set-TextCell $ws.Range("Y13")
$ws.Range("Y13").Value = "2023-09-30"
Set-TextCell $ws.Range("AA13")
$ws.Range("AA13").Value = "2023-09-30"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = "Ida Jansson"
$ws.Range("AX13").Value = "Ida Jansson"
